$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.710.12"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.388.29"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D5").Value = "505.29"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "132.89"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "2.390.93"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").Value = "4.69"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "2.811.76"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "56.647.59"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "21.72"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "2.382.54"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "10.20"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "4.06"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "310.02"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -4.64%  "
$ws.Range("D25").Value = "66.27"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "0.371"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "7.37"
$ws.Range("D30").Value = "175.43"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").Value = "5.88"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").Value = "3.81"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("D40").Value = "36.78"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("D42").Value = "1.44"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "132.61"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "4.83"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "0.568"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "246.95"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "0.0485"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").Value = "17.19"
$ws.Range("E51").Value = "  +7.34%  "
